$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some Price column values look like plain numbers (e.g. "1.002") but
# must be preserved as literal text (matching the source sheet which stores
# all Price/Volume cells as text). A leading apostrophe forces Excel to treat
# the assigned value as text instead of auto-converting/rounding it as a number.

$ws.Range("D2").Value = "25.863.60"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "1.639.08"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'215.95"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").Value = "'0.5075"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.2579"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "'0.06440"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "'4.287"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "1.864.18"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "1.636.20"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").Value = "'0.5641"
$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("D16").Value = "'63.29"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").Value = "0.0₅7606"
$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").Value = "25.877.63"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "'194.93"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'4.334"
$ws.Range("E21").Value = "  -2.83%  "

$ws.Range("D22").Value = "'9.910"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").Value = "'6.130"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'1.777"
$ws.Range("E25").Value = "  -6.34%  "

$ws.Range("D26").Value = "'0.1276"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("D27").Value = "'140.23"

$ws.Range("D28").Value = "'6.787"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").Value = "'15.48"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "'1.242"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").Value = "'0.04888"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").Value = "'3.298"
$ws.Range("E32").Value = "  +1.43%  "

$ws.Range("D33").Value = "'3.223"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").Value = "'1.559"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").Value = "'2.367"
$ws.Range("E35").Value = "  -0.36%  "

$ws.Range("D36").Value = "'0.9062"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").Value = "'2.579"
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").Value = "1.129.58"
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("D39").Value = "'0.5509"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").Value = "'0.01565"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").Value = "'0.9945"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").Value = "'0.8008"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "'97.63"
$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("D45").Value = "1.773.31"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("E46").Value = "  -6.06%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.56"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4445"
$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("D49").Value = "'7.716"
$ws.Range("E49").Value = "  +2.60%  "

$ws.Range("D50").Value = "'0.05054"
$ws.Range("E50").Value = "  -2.44%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.29%  "
